$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "... Senior Vice President of Technical Staff for Solutions and
#    Software at Cypress Semiconductor. " -> " welcome back to WICED
#    WiFI101 . "
# ------------------------------------------------------------------
$range = $d.Content
$old1 = " Senior Vice President of Technical Staff for Solutions and Software at Cypress Semiconductor. "
$found1 = $range.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $range.Delete()
    $range.InsertAfter(" welcome back to WICED WiFI101 . ")
}

# ------------------------------------------------------------------
# 2) Remove the existing "_GoBack" bookmark near the end of the
#    document -- it gets relocated to the "Things" paragraph below,
#    mirroring the last edit position Word itself would track.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 3) Append new sentences after "...a batch operation too." and drop
#    the "_GoBack" bookmark right after "...make your company's".
# ------------------------------------------------------------------
$range2 = $d.Content
$found2 = $range2.Find.Execute("a batch operation too.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $range2.Collapse(0)
    $insertText = "  And, all of this can also be done programmatically i.e. to make your company’s provisioning easier."
    $range2.InsertAfter($insertText)

    $marker = "company’s"
    $markerPos = $insertText.IndexOf($marker) + $marker.Length
    $bmPos = $range2.Start + $markerPos
    $bmRange = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# ------------------------------------------------------------------
# 4) Merge the three runs of the "resource ARN" sentence into one.
# ------------------------------------------------------------------
$range3 = $d.Content
$old3 = "For the resource ARN I'll just put in * and I'll check " + [char]34 + "Allow" + [char]34 + " so that it will be able to update any resource. Then I'll click on Create."
$found3 = $range3.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $range3.Delete()
    $range3.InsertAfter($old3)
}
